# Refreshes the crypto price ("D") and 1h-volume-change ("E") columns with
# the latest scrape, and swaps the ARBITRUM/Stacks rows (40-41) back to the
# order reported by the upstream feed.
#
# Price cells are plain text in the workbook (e.g. "357.80", "14.00"), so we
# force a Text number format before assigning any numeric-looking string -
# otherwise Excel auto-converts it to a Double and mangles the formatting
# (trailing zeros, thousand separators, floating point noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.296.88"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.931.37"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.80"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.74"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.09"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0872"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.63"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "3.387.66"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "2.922.01"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.990"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "52.281.00"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.60"
$ws.Range("E19").Value = "  +8.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  -1.54%  "
$ws.Range("D22").Value = "0.0₃0986"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.73"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.51"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.185"
$ws.Range("E26").Value = "  +7.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.71"
$ws.Range("E27").Value = "  +15.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.02"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.108"
$ws.Range("E30").Value = "  +9.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.52"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.59"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.20"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.22"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.29"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0445"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.36"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.00"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.99"
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.22"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.48"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  -5.59%  "
$ws.Range("D48").Value = "2.131.98"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.252"
$ws.Range("E49").Value = "  -5.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0347"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.931"
$ws.Range("E51").Value = "  -3.89%  "
